$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

for ($r = 3; $r -le 33; $r++) {
    $ws.Range("AJ$r").Formula = "=AH$r/AI$r"
}

$ws.Range("AL5").Select()
